$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new headers for columns I and J, reusing the same header formatting as H1
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data for columns I (I0) and J (IF), one pair per row from row 2 to row 33
$values = @(
    @(1,5),
    @(1,6),
    @(1,9),
    @(1,6),
    @(1,7),
    @(1,9),
    @(1,7),
    @(1,7),
    @(1,5),
    @(10,10),
    @(1,6),
    @(5,6),
    @(9,9),
    @(8,9),
    @(9,9),
    @(9,9),
    @(9,9),
    @(2,5),
    @(9,9),
    @(4,7),
    @(1,4),
    @(3,6),
    @(6,7),
    @(5,7),
    @(4,6),
    @(1,3),
    @(1,6),
    @(1,6),
    @(1,6),
    @(1,4),
    @(5,7),
    @(1,2)
)

for ($idx = 0; $idx -lt $values.Count; $idx++) {
    $row = $idx + 2
    $pair = $values[$idx]
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
}
